# Working on starve feeding section of methodology
# Adds a new "Material" / "Sourcing" table (rows 25-32) below the existing
# extrudability table on "Sheet1", with a blank spacer row (24) and a
# thick-bottom-border closing row (33).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New table data: Column A = component/material, Column B = sourcing.
$ws.Range("A25").Value = "Material"
$ws.Range("B25").Value = "Sourcing"

$ws.Range("A26").Value = "Device Body"
$ws.Range("B26").Value = "3D Printed (PLA)"

$ws.Range("A27").Value = "Ball Bearings"
$ws.Range("B27").Value = "Labs or 3D Printed (PLA)"

$ws.Range("A28").Value = "Gearbox"
$ws.Range("B28").Value = "Labs or 3D Printed (PLA)"

$ws.Range("A29").Value = "Drive Belt"
$ws.Range("B29").Value = "3D Printed (TPU)"

$ws.Range("A30").Value = "Ring Stand"
$ws.Range("B30").Value = "Labs"

$ws.Range("A31").Value = "Stepper Motor"
$ws.Range("B31").Value = "Labs"

$ws.Range("A32").Value = "UI Components"
$ws.Range("B32").Value = "Labs"

# Match styling of the existing tables: header row uses style 11/36 for the
# data rows; copy formatting from the analogous rows of the table above
# (rows 1-6 use style 24/8/9 header, 11/36/10 body, 12/13/14 closing border).
$ws.Range("A24:C24").Copy() | Out-Null
$ws.Range("A7:C7").Copy($ws.Range("A24:C24"))

$ws.Range("A25:C32").NumberFormat = "General"
for ($r = 25; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item(2, 1).Style
    $ws.Cells.Item($r, 2).Style = $ws.Cells.Item(2, 2).Style
    $ws.Cells.Item($r, 3).Style = $ws.Cells.Item(2, 3).Style
}

$ws.Range("A33:C33").Copy() | Out-Null
$ws.Range("A23:I23").Copy($ws.Range("A33:C33"))

# Autofit the first two columns to match the bestFit column widths.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Restore view state (scroll position, zoom, selection) to match final file.
$ws.Application.ActiveWindow.Zoom = 238
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A25:B32").Select()
